# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$changes = @{
    "F2"  = 8435
    "F3"  = 8048
    "F4"  = 133
    "F5"  = 196
    "F10" = 184
    "F13" = 158
    "F14" = 2125
    "F16" = 66
    "F19" = 138
    "F20" = 66
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $changes.Keys) {
        $ws.Range($addr).Value = $changes[$addr]
    }
}
